$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The experiment table's first measurement group (N=8, M=5) was dropped,
# so every remaining row moves up by one row-pair (data row + blank
# separator row). Remove rows 3:4 to perform that shift.
$ws.Rows("3:4").Delete()

# After the shift, row 7 holds the old "N=64 / ND" bold-styled row and
# row 9 holds the old "N=8, M=10" plain-styled row - but in the target
# layout row 7 becomes the new N=32 (M=10) plain data row and row 9
# becomes the N=64 (M=10) bold "ND" row. Swap their formatting (only)
# using a scratch row far outside the table as a temporary holder.
$ws.Range("A7:N7").Copy()
$ws.Range("A100:N100").PasteSpecial(-4122)

$ws.Range("A9:N9").Copy()
$ws.Range("A7:N7").PasteSpecial(-4122)

$ws.Range("A100:N100").Copy()
$ws.Range("A9:N9").PasteSpecial(-4122)

$ws.Range("A100:N100").Clear()

# Fill in the freshly measured M=10 data points (commit: "Cambiato M = 10").
$ws.Range("A3").Value = 8
$ws.Range("C3").Value = 10
$ws.Range("E3").Value = "212,26 MHz"
$ws.Range("G3").Value = "4,227 ns"
$ws.Range("I3").Value = "93 mW"
$ws.Range("K3").Value = "6 mW"
$ws.Range("M3").Value = 0.2
$ws.Range("N3").Value = 0.0091

$ws.Range("A5").Value = 16
$ws.Range("C5").Value = 10
$ws.Range("E5").Value = "202,75 Mhz"
$ws.Range("G5").Value = "4,634 ns"
$ws.Range("I5").Value = "95 mW"
$ws.Range("K5").Value = "10 mW"
$ws.Range("M5").Value = 0.36
$ws.Range("N5").Value = 0.0152

$ws.Range("A7").Value = 32
$ws.Range("C7").Value = 10
$ws.Range("E7").Value = "176,08 MHz"
$ws.Range("G7").Value = "5,338 ns"
$ws.Range("I7").Value = "95 mW"
$ws.Range("K7").Value = "17 mW"
$ws.Range("M7").Value = 0.68
$ws.Range("N7").Value = 0.0263

$ws.Range("A9").Value = 64
$ws.Range("C9").Value = 10
$ws.Range("E9").Value = "ND"
$ws.Range("G9").Value = "ND"
$ws.Range("I9").Value = "ND"
$ws.Range("K9").Value = "ND"
$ws.Range("M9").Value = 1.32
$ws.Range("N9").Value = "ND"

# Row 19 (N=32, M=100) had its frequency label corrected.
$ws.Range("E19").Value = "136,4 MHz"

# Update the saved selection to match the authored workbook.
$ws.Range("S4").Select()
